# Two more data measurements were taken for the "45.5 cm" BPA test sheet.
# New rows 17 and 18 are added to the data table (A:G), the derived
# formula columns (C, E, F) are extended to cover the two new rows, and
# column D ("620") is carried down through the remaining blank rows
# (19, 20 and the newly-created row 21) underneath the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("45.5 cm")

# --- New measurement row 17 ---
$ws.Range("A17").Value = 452
$ws.Range("B17").Value = 98
$ws.Range("D17").Value = 620
$ws.Range("G17").Value = 702

# --- New measurement row 18 ---
$ws.Range("A18").Value = 450
$ws.Range("B18").Value = 91
$ws.Range("D18").Value = 620
$ws.Range("G18").Value = 700

# --- Extend the derived-value formulas down through row 18 ---
# Column E: relative strain = 1 - A/$B$1 (was shared over E7:E16)
$ws.Range("E7:E18").FormulaR1C1 = "=1-RC[-4]/R1C2"
# Column F: normalized strain = E/$B$3 (was shared over F7:F16)
$ws.Range("F7:F18").FormulaR1C1 = "=RC[-1]/R3C2"
# Column C: force in Newtons = B*4.4482216 (was shared over C14:C16)
$ws.Range("C14:C18").FormulaR1C1 = "=RC[-1]*4.4482216"

# --- Carry the 620 value in column D down through the trailing rows ---
$ws.Range("D19").Value = 620
$ws.Range("D20").Value = 620
$ws.Range("D21").Value = 620

$excel.Calculate()

# --- Point the "BPA" scatter-chart series at the now-larger ranges so it
#     picks up the two new points too ---
$co = $ws.ChartObjects().Item(1)
$chart = $co.Chart
$series1 = $chart.SeriesCollection().Item(1)
$series1.XValues = $ws.Range("F6:F18")
$series1.Values = $ws.Range("C6:C18")

# --- The "45.5 cm" sheet (the sheet being actively worked on) becomes
#     the selected/active tab, moving away from "41.5 cm" ---
$ws.Activate() | Out-Null
$ws.Range("B19").Select() | Out-Null
